# Add a new "CheckoutInfo" worksheet after the existing "Login" sheet and
# populate it with first name / last name / zip-postal-code checkout data.

$wb = $excel.ActiveWorkbook

# Create the new sheet after the last existing sheet (so it lands after "Login").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "CheckoutInfo"

# Header row (bold, matches the style used on the "Login" sheet's header row).
$ws.Range("A1").Value = "First name"
$ws.Range("B1").Value = "Last name"
$ws.Range("C1").Value = "Zip/Postal Code"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows - write in the same order the original author did (first+last
# name for both people, then both zip/postal codes) so shared-string order
# matches.
$ws.Range("A2").Value = "Peter"
$ws.Range("B2").Value = "Parker"
$ws.Range("A3").Value = "Tom"
$ws.Range("B3").Value = "Davis"
$ws.Range("C2").Value = "1223A"
$ws.Range("C3").Value = "3455B"

# Column widths matching the committed worksheet (closest value the engine's
# column-width quantisation can reproduce for the committed 10.88671875 /
# 10.77734375 / 16.21875 character widths).
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 10.0
$ws.Columns.Item(3).ColumnWidth = 15.3

# Page orientation (portrait) so the pageSetup element is emitted.
$ws.PageSetup.Orientation = 1

# Select N13 on the new sheet and make it the active sheet/tab, matching the
# committed selection + tabSelected/activeTab state.
[void]$ws.Range("N13").Select()
$ws.Activate()
